# Update the Handback status report timestamps ("Generate Report for Handback")
$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the first data row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-16 09:05:44"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-16 09:05:38"
$wsZhCn.Range("K2").Value = "2016-08-16 09:05:55"

# de-de sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-16 09:05:44"
$wsDeDe.Range("K2").Value = "2016-08-16 09:06:05"
